$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.346.04'
$ws.Range("E2").Value = '  -0.66%  '
$ws.Range("D3").Value = '1.711.49'
$ws.Range("E3").Value = '  -0.58%  '
$ws.Range("E4").Value = '  +0.13%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '224.70'
$ws.Range("E5").Value = '  -0.53%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5294'
$ws.Range("E6").Value = '  -0.69%  '
$ws.Range("E7").Value = '  +0.14%  '
$ws.Range("B8").Value = 'Dogecoin'
$ws.Range("C8").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.06671'
$ws.Range("E8").Value = '  +1.27%  '
$ws.Range("B9").Value = 'Cardano'
$ws.Range("C9").Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.2664'
$ws.Range("E9").Value = '  -0.21%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '20.84'
$ws.Range("E10").Value = '  -3.70%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07699'
$ws.Range("E11").Value = '  -0.41%  '
$ws.Range("E12").Value = '  -2.47%  '
$ws.Range("D13").Value = '1.947.12'
$ws.Range("E13").Value = '  -0.60%  '
$ws.Range("D14").Value = '1.710.45'
$ws.Range("E14").Value = '  -0.65%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.5835'
$ws.Range("E15").Value = '  +0.10%  '
$ws.Range("D16").Value = '0.0₅8223'
$ws.Range("E16").Value = '  -0.68%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '68.05'
$ws.Range("E17").Value = '  +0.29%  '
$ws.Range("D18").Value = '27.370.99'
$ws.Range("E18").Value = '  -0.58%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '221.84'
$ws.Range("E19").Value = '  +0.88%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '1.004'
$ws.Range("E20").Value = '  +0.03%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.637'
$ws.Range("E21").Value = '  -1.98%  '
$ws.Range("E22").Value = '  -1.70%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.014'
$ws.Range("E23").Value = '  -0.96%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.006'
$ws.Range("E24").Value = '  +0.11%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '144.42'
$ws.Range("E25").Value = '  -1.52%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.695'
$ws.Range("E26").Value = '  -2.19%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.1204'
$ws.Range("E27").Value = '  -2.43%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.233'
$ws.Range("E28").Value = '  -2.35%  '
$ws.Range("E29").Value = '  -1.93%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.05328'
$ws.Range("E30").Value = '  -3.93%  '
$ws.Range("E31").Value = '  -0.84%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.461'
$ws.Range("E32").Value = '  -2.75%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.435'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.635'
$ws.Range("E34").Value = '  -1.60%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.874'
$ws.Range("E35").Value = '  +0.52%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.9521'
$ws.Range("E36").Value = '  -1.24%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.394'
$ws.Range("E37").Value = '  -0.98%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.5855'
$ws.Range("E38").Value = '  -1.85%  '
$ws.Range("D39").Value = '1.145.17'
$ws.Range("E39").Value = '  +8.39%  '
$ws.Range("E40").Value = '  -1.15%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.813'
$ws.Range("E42").Value = '  +0.14%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.8412'
$ws.Range("E43").Value = '  -1.51%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '101.33'
$ws.Range("E44").Value = '  -0.07%  '
$ws.Range("D45").Value = '1.854.17'
$ws.Range("E45").Value = '  -0.63%  '
$ws.Range("E46").Value = '  -0.77%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '57.79'
$ws.Range("E47").Value = '  -1.90%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.4549'
$ws.Range("E48").Value = '  +2.52%  '
$ws.Range("E49").Value = '  -0.34%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '8.060'
$ws.Range("E50").Value = '  -2.11%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.05223'
$ws.Range("E51").Value = '  -0.35%  '
